$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws "D2" '68.597.44'
Set-TextValue $ws "E2" '  -1.48%  '

# Row 3
Set-TextValue $ws "D3" '3.856.69'
Set-TextValue $ws "E3" '  -0.95%  '

# Row 4
Set-TextValue $ws "E4" '  -0.07%  '

# Row 5
Set-TextValue $ws "D5" '602.79'

# Row 6
Set-TextValue $ws "D6" '169.59'
Set-TextValue $ws "E6" '  -0.79%  '

# Row 7
Set-TextValue $ws "D7" '3.854.85'
Set-TextValue $ws "E7" '  -0.94%  '

# Row 8
Set-TextValue $ws "E8" '  +0.00%  '

# Row 9
Set-TextValue $ws "D9" '0.529'
Set-TextValue $ws "E9" '  -1.36%  '

# Row 10
Set-TextValue $ws "D10" '0.166'
Set-TextValue $ws "E10" '  -1.64%  '

# Row 11
Set-TextValue $ws "E11" '  +1.52%  '

# Row 12
Set-TextValue $ws "E12" '  -1.97%  '

# Row 13
Set-TextValue $ws "D13" '0.0000268'
Set-TextValue $ws "E13" '  +4.70%  '

# Row 14
Set-TextValue $ws "D14" '37.14'
Set-TextValue $ws "E14" '  -3.06%  '

# Row 15
Set-TextValue $ws "D15" '4.498.91'
Set-TextValue $ws "E15" '  -1.01%  '

# Row 16
Set-TextValue $ws "D16" '3.850.90'
Set-TextValue $ws "E16" '  -0.83%  '

# Row 17
Set-TextValue $ws "D17" '68.718.07'
Set-TextValue $ws "E17" '  -1.35%  '

# Row 18
Set-TextValue $ws "D18" '18.48'
Set-TextValue $ws "E18" '  -1.05%  '

# Row 19
Set-TextValue $ws "D19" '7.40'
Set-TextValue $ws "E19" '  -2.93%  '

# Row 20
Set-TextValue $ws "E20" '  -0.88%  '

# Row 21
Set-TextValue $ws "D21" '11.18'
Set-TextValue $ws "E21" '  +0.82%  '

# Row 22
Set-TextValue $ws "D22" '471.02'
Set-TextValue $ws "E22" '  -4.07%  '

# Row 23
Set-TextValue $ws "D23" '0.734'
Set-TextValue $ws "E23" '  -1.61%  '

# Row 24
Set-TextValue $ws "E24" '  -1.26%  '

# Row 25
Set-TextValue $ws "D25" '83.51'
Set-TextValue $ws "E25" '  -2.08%  '

# Row 26
Set-TextValue $ws "E26" '  -2.80%  '

# Row 27
Set-TextValue $ws "E27" '  -1.34%  '

# Row 28
Set-TextValue $ws "D28" '10.17'
Set-TextValue $ws "E28" '  +0.22%  '

# Row 29
Set-TextValue $ws "E29" '  +0.13%  '

# Row 30
Set-TextValue $ws "E30" '  -0.46%  '

# Row 31
Set-TextValue $ws "D31" '4.003.38'
Set-TextValue $ws "E31" '  -0.93%  '

# Row 32
Set-TextValue $ws "D32" '7.69'
Set-TextValue $ws "E32" '  -2.36%  '

# Row 33
Set-TextValue $ws "D33" '31.53'
Set-TextValue $ws "E33" '  -1.40%  '

# Row 34
Set-TextValue $ws "D34" '2.31'
Set-TextValue $ws "E34" '  -4.36%  '

# Row 35
Set-TextValue $ws "D35" '9.36'
Set-TextValue $ws "E35" '  -2.96%  '

# Row 36
Set-TextValue $ws "D36" '3.817.85'
Set-TextValue $ws "E36" '  -1.04%  '

# Row 37
Set-TextValue $ws "E37" '  -2.04%  '

# Row 38
Set-TextValue $ws "E38" '  +8.98%  '

# Row 39
Set-TextValue $ws "E39" '  -1.44%  '

# Row 40
Set-TextValue $ws "E40" '  -1.92%  '

# Row 41
Set-TextValue $ws "D41" '5.94'
Set-TextValue $ws "E41" '  -2.81%  '

# Row 42
Set-TextValue $ws "D42" '1.00'
Set-TextValue $ws "E42" '  +0.02%  '

# Row 43
Set-TextValue $ws "D43" '0.315'
Set-TextValue $ws "E43" '  -3.96%  '

# Row 44
Set-TextValue $ws "E44" '  -5.13%  '

# Row 45
Set-TextValue $ws "D45" '8.73'
Set-TextValue $ws "E45" '  +0.42%  '

# Row 46
$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextValue $ws "D46" '0.000295'
Set-TextValue $ws "E46" '  +6.95%  '

# Row 47
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws "D47" '417.60'
Set-TextValue $ws "E47" '  -3.73%  '

# Row 48
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws "D48" '1.00'
Set-TextValue $ws "E48" '  +0.00%  '

# Row 49
Set-TextValue $ws "E49" '  -1.96%  '

# Row 50
Set-TextValue $ws "E50" '  -1.60%  '

# Row 51
Set-TextValue $ws "D51" '141.83'
Set-TextValue $ws "E51" '  -1.08%  '
